# "Generate Report for Handback"
#
# The handback CI run failed to transform the file
# d01431d5-81bf-4829-a91d-a46d30f5e4fe.f91f55c1de8fd2bb2b45138b13dd753a0ee3dfa0
# for both the zh-cn and de-de locales, because the handback filename
# (ylwqoc3r.3hd) didn't match the handoff filename. Update the status
# report workbook accordingly:
#   1. Flip the "Status" for that file, on every sheet that reports it
#      (Overview, zh-cn, de-de), from "Ready for handoff" to
#      "Handback transform failed".
#   2. Record the failure reason in the "Error Detail" column (K) of the
#      per-locale sheets, row 3.

$wb = $excel.ActiveWorkbook

$statusFailed = "Handback transform failed"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusFailed
$wsOverview.Range("C3").Value = $statusFailed

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusFailed
$wsZhCn.Range("K3").Value = "Handback file name: ylwqoc3r.3hd is different with handoff file name: d01431d5-81bf-4829-a91d-a46d30f5e4fe.f91f55c1de8fd2bb2b45138b13dd753a0ee3dfa0.zh-cn."

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusFailed
$wsDeDe.Range("K3").Value = "Handback file name: ylwqoc3r.3hd is different with handoff file name: d01431d5-81bf-4829-a91d-a46d30f5e4fe.f91f55c1de8fd2bb2b45138b13dd753a0ee3dfa0.de-de."
